$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in new columns D, E, F, G for rows 2-7
$ws.Range("D2").Value = "F"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1

$ws.Range("D3").Value = "F"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 2

$ws.Range("D4").Value = "F"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 2
$ws.Range("G4").Value = 1

$ws.Range("D5").Value = "F"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2

$ws.Range("D6").Value = "F"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 3
$ws.Range("G6").Value = 3

$ws.Range("D7").Value = "F"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 4
$ws.Range("G7").Value = 4

# Update selection to F13
$ws.Range("F13").Select()
